$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.361.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.575.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.574.68"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.17"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.180.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.574.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.369.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.77"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.719.69"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.18"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.94"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +27.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.65"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.582.21"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.27"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.98"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.14"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.56"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0821"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.830"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.51"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.37%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.440.40"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.25%  "